{"js": "// Replace title, pros/cons list items, and the trailing bold/italic\n// SEO title + meta description paragraphs, per the commit diff.\n\nconst replacements = [\n  [\n    \"Play Horror Hotel Free: Review of Slot Game with High-Quality Graphics\",\n    \"Play Horror Hotel for Free - Review and Gameplay Analysis\",\n  ],\n  [\n    \"High-quality symbol design and graphics\",\n    \"Cluster system adds a unique twist to gameplay\",\n  ],\n  [\n    \"Fun and intriguing bonus game\",\n    \"High volatility for potential big wins\",\n  ],\n  [\n    \"Responsive design for gameplay on mobile devices\",\n    \"Well-designed symbols with impressive animations\",\n  ],\n  [\n    \"Above-average theoretical return to player percentage\",\n    \"Responsive design for seamless mobile play\",\n  ],\n  [\n    \"High volatility may be frustrating for some players\",\n    \"Winning clusters can be infrequent due to high volatility\",\n  ],\n  [\n    \"Autoplay and Turbo features may feel repetitive after extended gameplay\",\n    \"Bonus game can be difficult to trigger\",\n  ],\n  [\n    \"Read a detailed review of Horror Hotel, a high-quality online slot game with responsive design. Play for free and win big at this horror-themed slot.\",\n    \"Discover the features and gameplay of Horror Hotel slot. Play for free and test your luck!\",\n  ],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace title, pros/cons list items, and the trailing bold/italic\n# SEO title + meta description paragraphs, per the commit diff.\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute([ref]$findText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, $replaceText, 2) | Out-Null\n}\n\nReplace-AllText \"Play Horror Hotel Free: Review of Slot Game with High-Quality Graphics\" \"Play Horror Hotel for Free - Review and Gameplay Analysis\"\nReplace-AllText \"High-quality symbol design and graphics\" \"Cluster system adds a unique twist to gameplay\"\nReplace-AllText \"Fun and intriguing bonus game\" \"High volatility for potential big wins\"\nReplace-AllText \"Responsive design for gameplay on mobile devices\" \"Well-designed symbols with impressive animations\"\nReplace-AllText \"Above-average theoretical return to player percentage\" \"Responsive design for seamless mobile play\"\nReplace-AllText \"High volatility may be frustrating for some players\" \"Winning clusters can be infrequent due to high volatility\"\nReplace-AllText \"Autoplay and Turbo features may feel repetitive after extended gameplay\" \"Bonus game can be difficult to trigger\"\nReplace-AllText \"Read a detailed review of Horror Hotel, a high-quality online slot game with responsive design. Play for free and win big at this horror-themed slot.\" \"Discover the features and gameplay of Horror Hotel slot. Play for free and test your luck!\"\n"}
